# Added filtering options for the Component Analysis
#
# The "naive_err" forecast-error sheet lists, for each vintage date (column A),
# the error for forecast horizons Q0..Q9 (columns B..K). This change applies a
# horizon filter that truncates each vintage's row to only the quarters that
# are actually in scope for that vintage, clearing out the now out-of-range
# trailing error values while leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> first column to clear through column K (trailing run of cells
# to remove for that row), derived from the applied filter.
$clearRanges = @(
    "I2:K2",
    "G3:K3",
    "I4:K4",
    "G5:K5",
    "I6:K6",
    "G7:K7",
    "I8:K8",
    "G9:K9",
    "I10:K10",
    "G11:K11",
    "I12:K12",
    "K13:K13",
    "I14:K14",
    "K15:K15",
    "I16:K16",
    "K17:K17",
    "J18:K18",
    "I19:K19",
    "K21:K21",
    "J22:K22",
    "I23:K23",
    "K25:K25",
    "J26:K26",
    "I27:K27",
    "K29:K29",
    "J30:K30",
    "I31:K31",
    "K33:K33",
    "J34:K34",
    "I35:K35",
    "K37:K37",
    "J38:K38",
    "I39:J39"
)

foreach ($range in $clearRanges) {
    $ws.Range($range).ClearContents()
}
